$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Rodada 6" column header in G1, matching the style/format
# of the other round headers (e.g. F1: bold, centered, bordered).
$ws.Range("G1").Value = "Rodada 6"
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)  # xlPasteFormats

# Fill in the "Rodada 6" scores (column G) for each team. Teams that have
# no score for this round (the same rows already missing a "Rodada 5"/F
# value) are left blank, matching the source data.
$rodada6 = @{
  2  = 99.56005859375
  3  = 92.06005859375
  4  = 106.14990234375
  5  = 75.259765625
  6  = 104.2001953125
  7  = 78.919921875
  8  = 58.300048828125
  9  = 103.66015625
  10 = 88.759765625
  12 = 100.85986328125
  13 = 119.85009765625
  14 = 108.56005859375
  15 = 88.35986328125
  17 = 77.5
  19 = 82.64990234375
  20 = 78.56005859375
  21 = 86.06005859375
  22 = 106.06005859375
  24 = 48.840087890625
  25 = 106.14990234375
  26 = 110.35986328125
  28 = 95.31982421875
  29 = 78.85986328125
  30 = 98.06005859375
  31 = 86.9599609375
  32 = 87.759765625
  33 = 92.0498046875
}

foreach ($row in $rodada6.Keys) {
    $ws.Cells.Item($row, 7).Value = $rodada6[$row]
}
